$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $value, $refAddr)
    $cell = $ws.Range($addr)
    $ref = $ws.Range($refAddr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $ref.Style
}

Set-TextValue $ws "D2" "42.147.09" "C2"
Set-TextValue $ws "E2" "  -1.65%  " "C2"
Set-TextValue $ws "D3" "2.290.01" "C3"
Set-TextValue $ws "E3" "  -2.90%  " "C3"
Set-TextValue $ws "E4" "  +0.10%  " "C4"
Set-TextValue $ws "D5" "316.64" "C5"
Set-TextValue $ws "E5" "  -0.59%  " "C5"
Set-TextValue $ws "D6" "103.44" "C6"
Set-TextValue $ws "E6" "  -3.60%  " "C6"
Set-TextValue $ws "D7" "0.630" "C7"
Set-TextValue $ws "E7" "  -0.80%  " "C7"
Set-TextValue $ws "E8" "  +0.02%  " "C8"
Set-TextValue $ws "D9" "0.605" "C9"
Set-TextValue $ws "E9" "  -2.53%  " "C9"
Set-TextValue $ws "D10" "39.32" "C10"
Set-TextValue $ws "E10" "  -4.93%  " "C10"
Set-TextValue $ws "E11" "  -2.61%  " "C11"
Set-TextValue $ws "D12" "8.24" "C12"
Set-TextValue $ws "E12" "  -3.78%  " "C12"
Set-TextValue $ws "E13" "  -0.50%  " "C13"
Set-TextValue $ws "D14" "0.958" "C14"
Set-TextValue $ws "E14" "  -4.58%  " "C14"
Set-TextValue $ws "D15" "15.19" "C15"
Set-TextValue $ws "E15" "  -4.55%  " "C15"
Set-TextValue $ws "D16" "2.635.43" "C16"
Set-TextValue $ws "E16" "  -2.79%  " "C16"
Set-TextValue $ws "D17" "2.291.41" "C17"
Set-TextValue $ws "E17" "  -2.80%  " "C17"
Set-TextValue $ws "D18" "42.029.01" "C18"
Set-TextValue $ws "E18" "  -1.88%  " "C18"
Set-TextValue $ws "D19" "7.37" "C19"
Set-TextValue $ws "E19" "  -2.79%  " "C19"
Set-TextValue $ws "D20" "0.0000106" "C20"
Set-TextValue $ws "E20" "  -0.89%  " "C20"
Set-TextValue $ws "E21" "  -0.73%  " "C21"
Set-TextValue $ws "D22" "73.16" "C22"
Set-TextValue $ws "E22" "  -3.92%  " "C22"
Set-TextValue $ws "D23" "277.08" "C23"
Set-TextValue $ws "E23" "  +4.10%  " "C23"
Set-TextValue $ws "D24" "10.13" "C24"
Set-TextValue $ws "E24" "  +7.40%  " "C24"
Set-TextValue $ws "D25" "2.26" "C25"
Set-TextValue $ws "E25" "  -2.31%  " "C25"
Set-TextValue $ws "E26" "  +0.67%  " "C26"
Set-TextValue $ws "E27" "  -5.80%  " "C27"
Set-TextValue $ws "D28" "2.36" "C28"
Set-TextValue $ws "E28" "  +4.90%  " "C28"
Set-TextValue $ws "D29" "22.79" "C29"
Set-TextValue $ws "E29" "  -2.62%  " "C29"
Set-TextValue $ws "E30" "  -1.86%  " "C30"
Set-TextValue $ws "D31" "162.79" "C31"
Set-TextValue $ws "E31" "  -2.94%  " "C31"
Set-TextValue $ws "E32" "  -4.45%  " "C32"
Set-TextValue $ws "E33" "  -1.81%  " "C33"
Set-TextValue $ws "D34" "5.80" "C34"
Set-TextValue $ws "E34" "  -3.15%  " "C34"
Set-TextValue $ws "E35" "  +3.62%  " "C35"
Set-TextValue $ws "E36" "  -4.76%  " "C36"
Set-TextValue $ws "D37" "4.50" "C37"
Set-TextValue $ws "E37" "  -4.76%  " "C37"
Set-TextValue $ws "E38" "  -4.58%  " "C38"
Set-TextValue $ws "D39" "2.82" "C39"
Set-TextValue $ws "E39" "  +3.90%  " "C39"
Set-TextValue $ws "D40" "3.72" "C40"
Set-TextValue $ws "E40" "  -3.10%  " "C40"
Set-TextValue $ws "D41" "98.97" "C41"
Set-TextValue $ws "E41" "  -6.45%  " "C41"
Set-TextValue $ws "E42" "  -4.79%  " "C42"
Set-TextValue $ws "D43" "69.03" "C43"
Set-TextValue $ws "E43" "  -2.81%  " "C43"
Set-TextValue $ws "E44" "  +0.15%  " "C44"
Set-TextValue $ws "E45" "  -6.29%  " "C45"
Set-TextValue $ws "D46" "112.56" "C46"
Set-TextValue $ws "E46" "  -0.65%  " "C46"
Set-TextValue $ws "E47" "  -4.04%  " "C47"
Set-TextValue $ws "D48" "77.20" "C48"
Set-TextValue $ws "E48" "  +1.48%  " "C48"
Set-TextValue $ws "E49" "  -2.96%  " "C49"
Set-TextValue $ws "D50" "5.27" "C50"
Set-TextValue $ws "E50" "  -4.76%  " "C50"
Set-TextValue $ws "D51" "1.580.92" "C51"
Set-TextValue $ws "E51" "  +0.53%  " "C51"
